$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value='45.674.36'; ForceText=$false},
    @{Cell="E2"; Value='  -2.79%  '; ForceText=$false},
    @{Cell="D3"; Value='2.358.91'; ForceText=$false},
    @{Cell="E3"; Value='  +0.79%  '; ForceText=$false},
    @{Cell="E4"; Value='  +0.08%  '; ForceText=$false},
    @{Cell="D5"; Value='299.33'; ForceText=$true},
    @{Cell="E5"; Value='  -2.01%  '; ForceText=$false},
    @{Cell="D6"; Value='97.22'; ForceText=$true},
    @{Cell="E6"; Value='  +0.36%  '; ForceText=$false},
    @{Cell="D7"; Value='0.569'; ForceText=$true},
    @{Cell="E7"; Value='  -1.58%  '; ForceText=$false},
    @{Cell="E8"; Value='  +0.10%  '; ForceText=$false},
    @{Cell="D9"; Value='0.512'; ForceText=$true},
    @{Cell="E9"; Value='  -5.17%  '; ForceText=$false},
    @{Cell="D10"; Value='34.53'; ForceText=$true},
    @{Cell="E10"; Value='  -3.72%  '; ForceText=$false},
    @{Cell="D11"; Value='0.0795'; ForceText=$true},
    @{Cell="E11"; Value='  -2.30%  '; ForceText=$false},
    @{Cell="D12"; Value='7.08'; ForceText=$true},
    @{Cell="E12"; Value='  -4.97%  '; ForceText=$false},
    @{Cell="E13"; Value='  -1.60%  '; ForceText=$false},
    @{Cell="D14"; Value='2.714.86'; ForceText=$false},
    @{Cell="E14"; Value='  +0.77%  '; ForceText=$false},
    @{Cell="D15"; Value='2.384.74'; ForceText=$false},
    @{Cell="E15"; Value='  +2.14%  '; ForceText=$false},
    @{Cell="D16"; Value='13.65'; ForceText=$true},
    @{Cell="E16"; Value='  -3.88%  '; ForceText=$false},
    @{Cell="D17"; Value='0.805'; ForceText=$true},
    @{Cell="E17"; Value='  -3.77%  '; ForceText=$false},
    @{Cell="D18"; Value='45.412.16'; ForceText=$false},
    @{Cell="E18"; Value='  -3.03%  '; ForceText=$false},
    @{Cell="D19"; Value='0.0₃0969'; ForceText=$false},
    @{Cell="E19"; Value='  +1.88%  '; ForceText=$false},
    @{Cell="D20"; Value='12.51'; ForceText=$true},
    @{Cell="E20"; Value='  -8.31%  '; ForceText=$false},
    @{Cell="D21"; Value='5.95'; ForceText=$true},
    @{Cell="E21"; Value='  -3.71%  '; ForceText=$false},
    @{Cell="D22"; Value='65.84'; ForceText=$true},
    @{Cell="E22"; Value='  -2.89%  '; ForceText=$false},
    @{Cell="D23"; Value='242.88'; ForceText=$true},
    @{Cell="E23"; Value='  -4.58%  '; ForceText=$false},
    @{Cell="D24"; Value='2.80'; ForceText=$true},
    @{Cell="E24"; Value='  -6.14%  '; ForceText=$false},
    @{Cell="D25"; Value='1.00'; ForceText=$true},
    @{Cell="E25"; Value='  +0.08%  '; ForceText=$false},
    @{Cell="D26"; Value='1.88'; ForceText=$true},
    @{Cell="E26"; Value='  -5.81%  '; ForceText=$false},
    @{Cell="D27"; Value='40.28'; ForceText=$true},
    @{Cell="E27"; Value='  -5.05%  '; ForceText=$false},
    @{Cell="D28"; Value='2.22'; ForceText=$true},
    @{Cell="E28"; Value='  -1.42%  '; ForceText=$false},
    @{Cell="D29"; Value='9.62'; ForceText=$true},
    @{Cell="E29"; Value='  -2.88%  '; ForceText=$false},
    @{Cell="D30"; Value='20.67'; ForceText=$true},
    @{Cell="E30"; Value='  +1.47%  '; ForceText=$false},
    @{Cell="D31"; Value='3.65'; ForceText=$true},
    @{Cell="E31"; Value='  +15.31%  '; ForceText=$false},
    @{Cell="E32"; Value='  +4.79%  '; ForceText=$false},
    @{Cell="D33"; Value='144.86'; ForceText=$true},
    @{Cell="E33"; Value='  -1.14%  '; ForceText=$false},
    @{Cell="D34"; Value='5.36'; ForceText=$true},
    @{Cell="E34"; Value='  -7.68%  '; ForceText=$false},
    @{Cell="D35"; Value='0.0770'; ForceText=$true},
    @{Cell="E35"; Value='  -5.81%  '; ForceText=$false},
    @{Cell="D36"; Value='0.110'; ForceText=$true},
    @{Cell="E36"; Value='  -3.61%  '; ForceText=$false},
    @{Cell="E37"; Value='  -3.32%  '; ForceText=$false},
    @{Cell="D38"; Value='1.81'; ForceText=$true},
    @{Cell="E38"; Value='  -0.56%  '; ForceText=$false},
    @{Cell="D39"; Value='15.31'; ForceText=$true},
    @{Cell="E39"; Value='  +10.34%  '; ForceText=$false},
    @{Cell="D40"; Value='3.87'; ForceText=$true},
    @{Cell="E40"; Value='  -3.76%  '; ForceText=$false},
    @{Cell="D41"; Value='0.0296'; ForceText=$true},
    @{Cell="E41"; Value='  -5.05%  '; ForceText=$false},
    @{Cell="D42"; Value='3.14'; ForceText=$true},
    @{Cell="E42"; Value='  -8.08%  '; ForceText=$false},
    @{Cell="D43"; Value='1.00'; ForceText=$true},
    @{Cell="E43"; Value='  +0.16%  '; ForceText=$false},
    @{Cell="D44"; Value='1.861.42'; ForceText=$false},
    @{Cell="E44"; Value='  +2.60%  '; ForceText=$false},
    @{Cell="D45"; Value='90.51'; ForceText=$true},
    @{Cell="E45"; Value='  -3.02%  '; ForceText=$false},
    @{Cell="D46"; Value='1.77'; ForceText=$true},
    @{Cell="E46"; Value='  -10.67%  '; ForceText=$false},
    @{Cell="D47"; Value='0.183'; ForceText=$true},
    @{Cell="E47"; Value='  -5.98%  '; ForceText=$false},
    @{Cell="D48"; Value='69.43'; ForceText=$true},
    @{Cell="E48"; Value='  -7.67%  '; ForceText=$false},
    @{Cell="D49"; Value='2.586.66'; ForceText=$false},
    @{Cell="E49"; Value='  +0.65%  '; ForceText=$false},
    @{Cell="B50"; Value='Aave'; ForceText=$false},
    @{Cell="C50"; Value='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; ForceText=$false},
    @{Cell="D50"; Value='96.20'; ForceText=$true},
    @{Cell="E50"; Value='  -3.03%  '; ForceText=$false},
    @{Cell="B51"; Value='FraxShare'; ForceText=$false},
    @{Cell="C51"; Value='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; ForceText=$false},
    @{Cell="D51"; Value='8.00'; ForceText=$true},
    @{Cell="E51"; Value='  -0.96%  '; ForceText=$false}
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        $ws.Range($u.Cell).NumberFormat = "@"
    }
    $ws.Range($u.Cell).Value = $u.Value
}
